$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 06:51:52"
$wsZhCn.Range("H2").Value = "2016-03-17 06:52:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 06:51:59"
$wsDeDe.Range("H2").Value = "2016-03-17 06:52:44"
